$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '91.468.36'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.17%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.172.19'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +2.55%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '240.29'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.09%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '622.52'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.01%  '
$ws.Range('E7').Value = '  +4.67%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.372'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.71%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.999'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.08%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '3.168.99'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.50%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.753'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.57%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.205'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.56%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000249'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.18%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '35.48'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.04%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.53'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.73%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '91.421.82'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.51%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.748.04'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.36%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.165.83'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.89%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.74'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -3.00%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '15.53'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +10.58%  '
$ws.Range('B21').Value = 'Polkadot'
$ws.Range('C21').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.91'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +5.97%  '
$ws.Range('B22').Value = 'PEPE'
$ws.Range('C22').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.0000210'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -5.81%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '444.80'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.97%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.23'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +3.27%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '6.15'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +3.89%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '89.25'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.02%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.12'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.10%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '3.271.39'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E29').Value = '  +0.00%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.140'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +58.36%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.233'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +17.06%  '
$ws.Range('E32').Value = '  +8.04%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '9.38'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.13%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.167'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +8.58%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '7.92'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +9.43%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '26.58'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.93%  '
$ws.Range('B37').Value = 'Bittensor'
$ws.Range('C37').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '515.20'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.20%  '
$ws.Range('B38').Value = 'MantraDAO'
$ws.Range('C38').Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.05'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +18.98%  '
$ws.Range('B39').Value = 'PolygonEcosystemToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.467'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +14.91%  '
$ws.Range('B40').Value = 'PancakeSwap'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.93'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.36%  '
$ws.Range('B41').Value = 'Fetch.AI'
$ws.Range('C41').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.34'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +4.41%  '
$ws.Range('B42').Value = 'dogwifhat'
$ws.Range('C42').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.49'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -9.58%  '
$ws.Range('B43').Value = 'Binance-PegBSC-USD'
$ws.Range('C43').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.798'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -20.24%  '
$ws.Range('E44').Value = '  -0.02%  '
$ws.Range('E45').Value = '  +0.00%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.718'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +3.29%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.93'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.30%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '155.95'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.24%  '
$ws.Range('E49').Value = '  +3.58%  '
$ws.Range('B50').Value = 'VeChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0332'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +14.04%  '
$ws.Range('B51').Value = 'Filecoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '4.47'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.44%  '
